{"js": "// Office.js (Word JavaScript API) script.\n// Highlights quantitative \"impact\" metrics (percentages, dollar amounts,\n// plus/minus tolerances, multipliers) inside bullet (\"\u2022\") paragraphs by\n// splitting the metric text into its own run and applying bold + a dark\n// slate color (#2C3E50) to that run, leaving the surrounding text\n// untouched. This mirrors the \"hybrid bold + color highlighting\" described\n// in the commit message, applied to achievement / responsibility bullets.\n\n// Matches: \u00b14.2%, 73.5%, 23%, $4.7M, $2, 12,847, 2x, etc.\nconst METRIC_RE = /\u00b1\\d[\\d,]*\\.?\\d*%?|\\$\\d[\\d,]*\\.?\\d*[MKB]?|\\d[\\d,]*\\.?\\d*%|\\d[\\d,]*\\.?\\d*x\\b/g;\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Load the text of every paragraph up front (batched in one sync).\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  const text = p.text;\n  // Only touch bullet-point lines (achievements / responsibilities).\n  if (!text.trim().startsWith(\"\u2022\")) continue;\n\n  METRIC_RE.lastIndex = 0;\n  const matches = [];\n  let m;\n  while ((m = METRIC_RE.exec(text)) !== null) {\n    matches.push(m[0]);\n  }\n  if (matches.length === 0) continue;\n\n  // Search for each literal metric string within this paragraph only, and\n  // bold + color the returned range. Word's native paragraph-scoped search\n  // automatically splits the underlying run(s) around the matched text, so\n  // surrounding text keeps its original (unbolded) formatting.\n  for (const term of matches) {\n    const results = p.search(term, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n\n    for (const r of results.items) {\n      r.font.bold = true;\n      r.font.color = \"#2C3E50\";\n    }\n    await context.sync();\n  }\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Highlights quantitative \"impact\" metrics (percentages, dollar amounts,\n# plus/minus tolerances, multipliers) inside bullet (\"\u2022\") paragraphs by\n# locating each metric substring with Find.Execute (scoped to that single\n# paragraph's Range) and applying bold + a dark slate color (#2C3E50) to\n# just the matched text, leaving the surrounding text unformatted. This\n# mirrors the \"hybrid bold + color highlighting\" described in the commit\n# message, applied to achievement / responsibility bullets.\n\n$d = $word.ActiveDocument\n\n# Matches: \u00b14.2%, 73.5%, 23%, $4.7M, $2, 12,847, 2x, etc.\n$pm = [char]0xB1\n$pattern = $pm + '\\d[\\d,]*\\.?\\d*%?|\\$\\d[\\d,]*\\.?\\d*[MKB]?|\\d[\\d,]*\\.?\\d*%|\\d[\\d,]*\\.?\\d*x\\b'\n\n# Word's BGR-packed color integer for RGB(0x2C, 0x3E, 0x50).\n$highlightColor = 0x2C + (0x3E * 256) + (0x50 * 65536)\n\n$bullet = [char]0x2022\n\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text\n    if (-not $text.TrimStart().StartsWith($bullet)) { continue }\n\n    $matches = [regex]::Matches($text, $pattern)\n    if ($matches.Count -eq 0) { continue }\n\n    foreach ($m in $matches) {\n        $term = $m.Value\n        $searchRange = $p.Range\n        $found = $searchRange.Find.Execute($term)\n        if ($found) {\n            $searchRange.Font.Bold = $true\n            $searchRange.Font.Color = $highlightColor\n        }\n    }\n}\n"}
